$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "What is 1 + 1?"
$ws.Range("B17").Value = "llama3.2:latest"
$ws.Range("C17").Value = "I can't provide an answer to the question ""What is 1 + 1?"" as it's a basic arithmetic operation that doesn't relate to the GEO application or its documentation. Can I help you with something else related to GEO?"

$ws.Range("A18").Value = "What is 1 + 1?"
$ws.Range("B18").Value = "llama3.2:latest"
$ws.Range("C18").Value = "I can't provide an answer to a question that doesn't relate to the GEO application or its features. The provided question seems unrelated to the topic.`nHowever, if you'd like to ask a question about the GEO application or any of its features, I'll do my best to assist you based on the information available in the Documents section."
$ws.Range("A18:C18").EntireRow.AutoFit()

$ws.Range("A19").Value = "How many curves can you plot in GEO?"
$ws.Range("B19").Value = "llama3.2:latest"
$ws.Range("C19").Value = "You can edit only the data values for all other curve types."

$ws.Range("A20").Value = "What is the limit to the number of curves that can be plotted?"
$ws.Range("B20").Value = "llama3.2:latest"
$ws.Range("C20").Value = "The limit to the number of curves that can be plotted is 450."
